# Update "想去人数" (want-to-go count) values for two events that appear
# on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览": row 3 -> F3 120 -> 122, row 4 -> F4 658 -> 663
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 122
$wsExhibit.Range("F4").Value = 663

# Sheet "全部类型": row 4 -> F4 120 -> 122, row 5 -> F5 658 -> 663
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 122
$wsAll.Range("F5").Value = 663
